# Fix 5.4.3 lab sheet formulas: replace SQRT(num)/SQRT(denom) with num/SQRT(denom2)
# using the correct "raw" denominator column for each block.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Block A-D (run 1)
$ws.Range("D3").Formula = "=C3/SQRT(B3)"
$ws.Range("D4:D5").Formula = "=C4/SQRT(B4)"

# Block E-H (run 2)
$ws.Range("H3").Formula = "=G3/SQRT(F3)"
$ws.Range("H4:H5").Formula = "=G4/SQRT(F4)"

# Block J-M (run 3)
$ws.Range("M3").Formula = "=L3/SQRT(K3)"
$ws.Range("M4:M6").Formula = "=L4/SQRT(K4)"

# Block N-Q (run 4)
$ws.Range("Q3").Formula = "=P3/SQRT(O3)"
$ws.Range("Q4:Q6").Formula = "=P4/SQRT(O4)"

# Block S-W (run 5)
$ws.Range("W3").Formula = "=V3/SQRT(U3)"
$ws.Range("W4:W8").Formula = "=V4/SQRT(U4)"

# Sheet view changes: scroll position and selection
$ws.Application.ActiveWindow.ScrollColumn = 19
$ws.Range("Z3").Select()
